$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "invalid email" test case (row 6) ---
# Title: drop the old wording, now reads "...validation mail"
$ws.Range("C6").Value = "Validating Registaration `nvalidation mail"

# Data: fix the typo'd email address (missing "@")
$ws.Range("D6").Value = "* First Name = ""Nehal""`n*LastName = ""Srour""`n* Phone = ""01028374659""`n* Email = ""noor@gmail.com""`n*Password = ""Noor%66""`n* Confirm Password ""Noor%66"""

# expected: append the validation-mail note
$ws.Range("G6").Value = "Error In the Registration,`n*Error Message ""Invalid Email"" is displayed`n* a validation mail is sent"

# Actual: append the "no validation mail" note
$ws.Range("H6").Value = "The Registration is done successfully`n*The Home Page is displayed`n* No Validation mail "

# --- Row 7 ("already registered email" test case) is removed ---
# Id / Related TC number columns go away completely (no leftover cell/style)
$ws.Range("A7:B7").Clear()
# The rest of the row is wiped of content but keeps its formatting
$ws.Range("C7:H7").ClearContents()

# --- View state: zoom out and move the selected cell ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 55
$ws.Range("D11").Select() | Out-Null
